$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00506879073135409
$ws.Range("C2").Value = 0.137581462708182
$ws.Range("D2").Value = 0.141926140477915
$ws.Range("E2").Value = 0.0883417813178856
$ws.Range("F2").Value = 0.00724112961622013
$ws.Range("G2").Value = 0.000724112961622013
$ws.Range("H2").Value = 0.86676321506155
$ws.Range("I2").Value = 0.00724112961622013
$ws.Range("J2").Value = 0.0057929036929761
$ws.Range("K2").Value = 0.858073859522085
$ws.Range("L2").Value = 0.0137581462708182
$ws.Range("M2").Value = 0.0246198406951484
$ws.Range("N2").Value = 0.000724112961622013
$ws.Range("O2").Value = 0.00868935553946416
$ws.Range("P2").Value = 0.807385952208545
$ws.Range("Q2").Value = 0.00434467776973208
$ws.Range("S2").Value = 0.997827661115134
$ws.Range("U2").Value = 0.441708906589428
$ws.Range("V2").Value = 0.0209992758870384
$ws.Range("W2").Value = 0.0926864590876177
$ws.Range("X2").Value = 0.00724112961622013
$ws.Range("B3").Value = 0.935553946415641
$ws.Range("C3").Value = 0.847936278059377
$ws.Range("D3").Value = 0.0115858073859522
$ws.Range("E3").Value = 0.865314989138306
$ws.Range("F3").Value = 0.973207820419986
$ws.Range("G3").Value = 0.00506879073135409
$ws.Range("H3").Value = 0.0159304851556843
$ws.Range("I3").Value = 0.0057929036929761
$ws.Range("J3").Value = 0.987690079652426
$ws.Range("K3").Value = 0.0608254887762491
$ws.Range("L3").Value = 0.0217233888486604
$ws.Range("M3").Value = 0.0137581462708182
$ws.Range("N3").Value = 0.174511223750905
$ws.Range("O3").Value = 0.97827661115134
$ws.Range("P3").Value = 0.0130340333091962
$ws.Range("Q3").Value = 0.0325850832729906
$ws.Range("R3").Value = 0.98913830557567
$ws.Range("T3").Value = 0.942795076031861
$ws.Range("U3").Value = 0.0231716147719044
$ws.Range("V3").Value = 0.0311368573497466
$ws.Range("W3").Value = 0.0188269370021723
$ws.Range("X3").Value = 0.0188269370021723
$ws.Range("B4").Value = 0.0188269370021723
$ws.Range("C4").Value = 0.00724112961622013
$ws.Range("D4").Value = 0.770456191165822
$ws.Range("E4").Value = 0.0376538740043447
$ws.Range("F4").Value = 0.0115858073859522
$ws.Range("G4").Value = 0.0057929036929761
$ws.Range("H4").Value = 0.11151339608979
$ws.Range("I4").Value = 0.986241853729182
$ws.Range("J4").Value = 0.000724112961622013
$ws.Range("K4").Value = 0.0709630702389573
$ws.Range("L4").Value = 0.0115858073859522
$ws.Range("M4").Value = 0.0101375814627082
$ws.Range("N4").Value = 0.0057929036929761
$ws.Range("O4").Value = 0.00651701665459812
$ws.Range("P4").Value = 0.178131788559015
$ws.Range("Q4").Value = 0.0173787110789283
$ws.Range("S4").Value = 0.00217233888486604
$ws.Range("T4").Value = 0.0173787110789283
$ws.Range("U4").Value = 0.529326574945692
$ws.Range("V4").Value = 0.0209992758870384
$ws.Range("W4").Value = 0.882693700217234
$ws.Range("X4").Value = 0.973207820419986
$ws.Range("B5").Value = 0.0405503258508327
$ws.Range("C5").Value = 0.00724112961622013
$ws.Range("D5").Value = 0.0745836350470673
$ws.Range("E5").Value = 0.00868935553946416
$ws.Range("F5").Value = 0.00796524257784214
$ws.Range("G5").Value = 0.988414192614048
$ws.Range("H5").Value = 0.00506879073135409
$ws.Range("I5").Value = 0.000724112961622013
$ws.Range("J5").Value = 0.00506879073135409
$ws.Range("K5").Value = 0.0101375814627082
$ws.Range("L5").Value = 0.952932657494569
$ws.Range("M5").Value = 0.951484431571325
$ws.Range("N5").Value = 0.818247646632875
$ws.Range("O5").Value = 0.00651701665459812
$ws.Range("P5").Value = 0.00144822592324403
$ws.Range("Q5").Value = 0.945691527878349
$ws.Range("R5").Value = 0.0108616944243302
$ws.Range("T5").Value = 0.0398262128892107
$ws.Range("U5").Value = 0.00506879073135409
$ws.Range("V5").Value = 0.926864590876177
$ws.Range("W5").Value = 0.0057929036929761
